# Auto-generated edit script applying the Phantom_Profits.xlsx diff
# Updates static price/profit data cells across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 5408.909
$ws.Range("I113").Value = 3799.2856
$ws.Range("J113").Value = 8225.75
$ws.Range("K113").Value = 3799.2856
$ws.Range("L113").Value = 8225.75
$ws.Range("M113").Value = -545.2856000000002
$ws.Range("N113").Value = -14733.75
# Row 116
$ws.Range("H116").Value = 4567.6
$ws.Range("I116").Value = 4656.75
$ws.Range("K116").Value = 4656.75
$ws.Range("M116").Value = -1214.75
# Row 137
$ws.Range("H137").Value = 2058.7778
$ws.Range("I137").Value = 2010.3334
$ws.Range("K137").Value = 6031.0002
$ws.Range("M137").Value = -3481.0002

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 7510.6
$ws.Range("I61").Value = 7888.25
$ws.Range("K61").Value = 7888.25
$ws.Range("M61").Value = -7676.25
# Row 132
$ws.Range("H132").Value = 4187.375
$ws.Range("I132").Value = 7666.6665
$ws.Range("J132").Value = 2099.8
$ws.Range("K132").Value = 22999.9995
$ws.Range("L132").Value = 6299.400000000001
$ws.Range("M132").Value = -20469.9995
$ws.Range("N132").Value = -11359.4
# Row 136
$ws.Range("H136").Value = 7510.6
$ws.Range("I136").Value = 7888.25
$ws.Range("K136").Value = 23664.75
$ws.Range("M136").Value = -21114.75

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3961.25
$ws.Range("I105").Value = 2659.4
$ws.Range("K105").Value = 2659.4
$ws.Range("M105").Value = -912.4000000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3759.9
# Row 34
$ws.Range("H34").Value = 3759.9
# Row 107
$ws.Range("H107").Value = 788.125
$ws.Range("I107").Value = 815.1429000000001
$ws.Range("K107").Value = 815.1429000000001
$ws.Range("M107").Value = 1104.8571

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 5580035
$ws.Range("I4").Value = 3695345
$ws.Range("J4").Value = 15003484
$ws.Range("K4").Value = 11086035
$ws.Range("L4").Value = 45010452
$ws.Range("M4").Value = -11085923
$ws.Range("N4").Value = -45010676
# Row 9
$ws.Range("H9").Value = 9312.5625
$ws.Range("J9").Value = 9312.5625
$ws.Range("L9").Value = 27937.6875
$ws.Range("N9").Value = -28385.6875
# Row 10
$ws.Range("H10").Value = 870
$ws.Range("I10").Value = 272.5
$ws.Range("J10").Value = 1666.6666
$ws.Range("K10").Value = 817.5
$ws.Range("L10").Value = 4999.9998
$ws.Range("M10").Value = -678.5
$ws.Range("N10").Value = -5277.9998
# Row 11
$ws.Range("H11").Value = 2073.5715
$ws.Range("I11").Value = 161
$ws.Range("J11").Value = 2392.3333
$ws.Range("K11").Value = 483
$ws.Range("L11").Value = 7176.999899999999
$ws.Range("M11").Value = -343
$ws.Range("N11").Value = -7456.999899999999
# Row 12
$ws.Range("H12").Value = 187.75
$ws.Range("I12").Value = 205.33333
$ws.Range("J12").Value = 135
$ws.Range("K12").Value = 615.99999
$ws.Range("L12").Value = 405
$ws.Range("M12").Value = -442.99999
$ws.Range("N12").Value = -751
# Row 13
$ws.Range("H13").Value = 7000
$ws.Range("I13").Value = 50
$ws.Range("J13").Value = 10475
$ws.Range("K13").Value = 150
$ws.Range("L13").Value = 31425
$ws.Range("M13").Value = 18
$ws.Range("N13").Value = -31761
# Row 14
$ws.Range("H14").Value = 1537.2142
$ws.Range("I14").Value = 1537.2142
$ws.Range("K14").Value = 4611.642599999999
$ws.Range("M14").Value = -4438.642599999999
# Row 15
$ws.Range("H15").Value = 152.5
$ws.Range("I15").Value = 201
$ws.Range("K15").Value = 603
$ws.Range("M15").Value = -463
# Row 16
$ws.Range("H16").Value = 560
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
# Row 17
$ws.Range("H17").Value = 916
# Row 70
$ws.Range("H70").Value = 17000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 17000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 51000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -51630
# Row 73
$ws.Range("H73").Value = 17000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 17000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 51000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -53184
# Row 80
$ws.Range("H80").Value = 4833.3335
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4833.3335
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 14500.0005
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -16372.0005
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
# Row 83
$ws.Range("H83").Value = 4833.3335
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4833.3335
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 43500.0015
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -52860.0015
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
# Row 98
$ws.Range("H98").Value = 357.14285
$ws.Range("I98").Value = 278.33334
$ws.Range("K98").Value = 835.0000200000001
$ws.Range("M98").Value = 662.9999799999999
# Row 112
$ws.Range("H112").Value = 4728.8335
$ws.Range("I112").Value = 1644
$ws.Range("J112").Value = 6271.25
$ws.Range("K112").Value = 4932
$ws.Range("L112").Value = 18813.75
$ws.Range("M112").Value = -3824
$ws.Range("N112").Value = -21029.75
# Row 113
$ws.Range("H113").Value = 1260.9166
$ws.Range("I113").Value = 536.5
$ws.Range("J113").Value = 1623.125
$ws.Range("K113").Value = 1609.5
$ws.Range("L113").Value = 4869.375
$ws.Range("M113").Value = 560.5
$ws.Range("N113").Value = -9209.375
# Row 116
$ws.Range("H116").Value = 4793.3335
$ws.Range("I116").Value = 4793.3335
$ws.Range("K116").Value = 14380.0005
$ws.Range("M116").Value = -10938.0005
# Row 119
$ws.Range("H119").Value = 509.33334
$ws.Range("I119").Value = 509.33334
$ws.Range("K119").Value = 1528.00002
$ws.Range("M119").Value = 3309.99998
# Row 120
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2450
$ws.Range("J80").Value = 2450
$ws.Range("L80").Value = 2450
$ws.Range("N80").Value = -4446
# Row 83
$ws.Range("H83").Value = 2450
$ws.Range("J83").Value = 2450
$ws.Range("L83").Value = 12250
$ws.Range("N83").Value = -22234
# Row 132
$ws.Range("H132").Value = 3107.074
$ws.Range("I132").Value = 3217.524
$ws.Range("J132").Value = 2720.5
$ws.Range("K132").Value = 9652.572
$ws.Range("L132").Value = 8161.5
$ws.Range("M132").Value = -7122.572
$ws.Range("N132").Value = -13221.5

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 7408.8
$ws.Range("I68").Value = 6761
$ws.Range("K68").Value = 6761
$ws.Range("M68").Value = -6012
# Row 71
$ws.Range("H71").Value = 7408.8
$ws.Range("I71").Value = 6761
$ws.Range("K71").Value = 33805
$ws.Range("M71").Value = -30061
# Row 106
$ws.Range("H106").Value = 27357.4
$ws.Range("J106").Value = 27357.4
$ws.Range("L106").Value = 27357.4
$ws.Range("N106").Value = -29881.4
# Row 136
$ws.Range("H136").Value = 31252132
$ws.Range("I136").Value = 2562.818
$ws.Range("J136").Value = 100001180
$ws.Range("K136").Value = 7688.454000000001
$ws.Range("L136").Value = 300003540
$ws.Range("M136").Value = -5138.454000000001
$ws.Range("N136").Value = -300008640

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1418.75
$ws.Range("I132").Value = 1418.75
$ws.Range("K132").Value = 4256.25
$ws.Range("M132").Value = -1726.25
